$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.699.11"
$ws.Range("E2").Value = "  +1.19%  "

$ws.Range("D3").Value = "3.841.73"
$ws.Range("E3").Value = "  +0.50%  "

$ws.Range("E4").Value = "  -0.39%  "

$ws.Range("D5").Value = "'448.12"
$ws.Range("E5").Value = "  +6.34%  "

$ws.Range("D6").Value = "'147.53"
$ws.Range("E6").Value = "  +15.48%  "

$ws.Range("D7").Value = "'0.627"
$ws.Range("E7").Value = "  +4.21%  "

$ws.Range("E8").Value = "  -0.18%  "

$ws.Range("D9").Value = "'0.750"
$ws.Range("E9").Value = "  +4.84%  "

$ws.Range("D10").Value = "'0.158"
$ws.Range("E10").Value = "  -2.32%  "

$ws.Range("D11").Value = "'0.0000325"
$ws.Range("E11").Value = "  -6.53%  "

$ws.Range("D12").Value = "'44.43"
$ws.Range("E12").Value = "  +11.62%  "

$ws.Range("D13").Value = "'10.49"
$ws.Range("E13").Value = "  +5.36%  "

$ws.Range("D14").Value = "4.447.26"
$ws.Range("E14").Value = "  -0.19%  "

$ws.Range("D15").Value = "'14.90"
$ws.Range("E15").Value = "  -8.41%  "

$ws.Range("D16").Value = "3.857.19"
$ws.Range("E16").Value = "  +1.02%  "

$ws.Range("E17").Value = "  -0.10%  "

$ws.Range("D18").Value = "'20.15"
$ws.Range("E18").Value = "  +4.09%  "

$ws.Range("E19").Value = "  +7.85%  "

$ws.Range("D20").Value = "67.738.27"
$ws.Range("E20").Value = "  +0.98%  "

$ws.Range("D21").Value = "'423.86"
$ws.Range("E21").Value = "  +4.66%  "

$ws.Range("D22").Value = "'14.76"
$ws.Range("E22").Value = "  +4.66%  "

$ws.Range("D23").Value = "'3.28"
$ws.Range("E23").Value = "  +10.39%  "

$ws.Range("D24").Value = "'86.84"
$ws.Range("E24").Value = "  +3.96%  "

$ws.Range("B25").Value = "EthereumClassic"
$ws.Range("C25").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D25").Value = "'37.71"
$ws.Range("E25").Value = "  +3.05%  "

$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "'3.47"
$ws.Range("E26").Value = "  +9.04%  "

$ws.Range("D27").Value = "'5.55"
$ws.Range("E27").Value = "  -5.64%  "

$ws.Range("D28").Value = "'9.88"
$ws.Range("E28").Value = "  +6.15%  "

$ws.Range("E29").Value = "  +23.09%  "

$ws.Range("D30").Value = "'731.96"
$ws.Range("E30").Value = "  +1.98%  "

$ws.Range("D31").Value = "'13.88"
$ws.Range("E31").Value = "  +13.21%  "

$ws.Range("E32").Value = "  +11.64%  "

$ws.Range("D33").Value = "'2.76"
$ws.Range("E33").Value = "  -0.22%  "

$ws.Range("D34").Value = "'43.98"
$ws.Range("E34").Value = "  +17.65%  "

$ws.Range("E35").Value = "  +7.53%  "

$ws.Range("D36").Value = "'56.82"
$ws.Range("E36").Value = "  +3.66%  "

$ws.Range("D37").Value = "'5.56"
$ws.Range("E37").Value = "  +22.58%  "

$ws.Range("E38").Value = "  +0.11%  "

$ws.Range("D39").Value = "'0.0482"
$ws.Range("E39").Value = "  +7.15%  "

$ws.Range("D40").Value = "'2.94"
$ws.Range("E40").Value = "  +1.46%  "

$ws.Range("D41").Value = "0.0₃0689"
$ws.Range("E41").Value = "  -9.40%  "

$ws.Range("E42").Value = "  +5.46%  "

$ws.Range("D43").Value = "'0.337"
$ws.Range("E43").Value = "  +16.27%  "

$ws.Range("E44").Value = "  -0.52%  "

$ws.Range("E45").Value = "  +2.77%  "

$ws.Range("B46").Value = "LidoDAOToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D46").Value = "'3.41"
$ws.Range("E46").Value = "  +3.23%  "

$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "'2.49"
$ws.Range("E47").Value = "  +18.70%  "

$ws.Range("E48").Value = "  +5.48%  "

$ws.Range("D49").Value = "'145.78"
$ws.Range("E49").Value = "  +1.55%  "

$ws.Range("D50").Value = "'2.92"
$ws.Range("E50").Value = "  +6.45%  "

$ws.Range("D51").Value = "'2.66"
$ws.Range("E51").Value = "  +6.16%  "
